$wb = $excel.ActiveWorkbook

# --- "Output" sheet: just move the selection cursor (no data changes) ---
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Range("O9").Select()

# --- "Login" sheet: the Data-Driven columns added to the right of Username/password ---
$wsLogin = $wb.Worksheets.Item("Login")

# Column headers (row 1) and values (row 2), entered in the same left-to-right /
# row-by-row order the original author used so the shared-string table lines up.
$wsLogin.Range("C1").Value = "customerNumber"
$wsLogin.Range("D1").Value = "coType"
$wsLogin.Range("C2").Value = "US00025065"
$wsLogin.Range("D2").Value = "USA"
$wsLogin.Range("E1").Value = "PO"
$wsLogin.Range("E2").Value = "US-Gear-06"
$wsLogin.Range("F1").Value = "ItemCode"
$wsLogin.Range("F2").Value = "TB7SX6CC"

# Highlight the "Username" header cell in red font.
$wsLogin.Range("A1").Font.Color = 255

# Best-fit the first two (original) columns.
$wsLogin.Columns.Item(1).ColumnWidth = 9.14
$wsLogin.Columns.Item(2).ColumnWidth = 10.42

# Print orientation for this sheet.
$wsLogin.PageSetup.Orientation = 1

# Leave the cursor on the Login tab (which stays the active sheet), matching
# the saved selection/tab state from the workbook.
$wsLogin.Range("F6").Select()
